# Updated cryptos list on Sat Feb 17 07:59:56 UTC 2024 with GitHub Actions
# Refresh of the Price (D) / Volume(1h) (E) columns from the latest feed,
# plus a ranking swap: Monero (rank 43) now edges out EnergySwap (rank 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D35 ("0.0850") and D47 ("3.30") carry a significant trailing zero that
# Excel's automatic type-detection would otherwise strip when parsing the
# assigned string as a Number. Pin those two cells to Text first so the
# literal source-feed string round-trips unchanged.
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range('D2').Value = '51.930.20'
$ws.Range('E2').Value = '  +0.33%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '2.796.99'
$ws.Range('E3').Value = '  -0.73%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  +0.00%  '

# Row 5 - BNB
$ws.Range('D5').Value = '359.99'
$ws.Range('E5').Value = '  +1.25%  '

# Row 6 - Solana
$ws.Range('D6').Value = '110.33'
$ws.Range('E6').Value = '  -1.11%  '

# Row 7 - XRP
$ws.Range('D7').Value = '0.561'
$ws.Range('E7').Value = '  -0.77%  '

# Row 8 - USDC
$ws.Range('E8').Value = '  -0.03%  '

# Row 9 - Cardano
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').Value = '  -0.90%  '

# Row 10 - Avalanche
$ws.Range('D10').Value = '40.27'
$ws.Range('E10').Value = '  -1.16%  '

# Row 11 - TRON
$ws.Range('E11').Value = '  +2.37%  '

# Row 12 - Dogecoin
$ws.Range('D12').Value = '0.0852'
$ws.Range('E12').Value = '  -0.03%  '

# Row 13 - Chainlink
$ws.Range('E13').Value = '  -1.52%  '

# Row 14 - Polkadot
$ws.Range('E14').Value = '  -0.83%  '

# Row 15 - Wrapped liquid staked Ether 2.0
$ws.Range('D15').Value = '3.234.34'
$ws.Range('E15').Value = '  -0.76%  '

# Row 16 - Wrapped Ether
$ws.Range('D16').Value = '2.793.26'
$ws.Range('E16').Value = '  -0.54%  '

# Row 17 - Polygon
$ws.Range('D17').Value = '0.948'
$ws.Range('E17').Value = '  +3.48%  '

# Row 18 - Wrapped BTC
$ws.Range('D18').Value = '51.897.33'
$ws.Range('E18').Value = '  +0.39%  '

# Row 19 - Uniswap
$ws.Range('D19').Value = '7.48'
$ws.Range('E19').Value = '  -0.59%  '

# Row 20 - ImmutableX
$ws.Range('E20').Value = '  -1.47%  '

# Row 21 - Internet Computer (DFINITY)
$ws.Range('D21').Value = '13.37'
$ws.Range('E21').Value = '  +0.36%  '

# Row 22 - Shiba Inu
$ws.Range('E22').Value = '  -0.90%  '

# Row 23 - Litecoin
$ws.Range('D23').Value = '70.38'
$ws.Range('E23').Value = '  +0.98%  '

# Row 24 - Bitcoin Cash
$ws.Range('D24').Value = '270.63'
$ws.Range('E24').Value = '  +1.28%  '

# Row 25 - PancakeSwap
$ws.Range('E25').Value = '  -0.65%  '

# Row 26 - Ethereum Classic
$ws.Range('D26').Value = '26.56'
$ws.Range('E26').Value = '  -1.35%  '

# Row 27 - Dai
$ws.Range('E27').Value = '  +0.03%  '

# Row 28 - Kaspa
$ws.Range('E28').Value = '  +19.07%  '

# Row 29 - Cosmos
$ws.Range('D29').Value = '10.29'
$ws.Range('E29').Value = '  +0.27%  '

# Row 30 - Toncoin
$ws.Range('E30').Value = '  -3.93%  '

# Row 31 - Filecoin
$ws.Range('D31').Value = '6.25'
$ws.Range('E31').Value = '  +6.42%  '

# Row 32 - Injective Protocol
$ws.Range('D32').Value = '35.11'
$ws.Range('E32').Value = '  +2.26%  '

# Row 33 - OKB
$ws.Range('D33').Value = '52.16'
$ws.Range('E33').Value = '  -0.24%  '

# Row 34 - VeChain
$ws.Range('D34').Value = '0.0465'
$ws.Range('E34').Value = '  -3.00%  '

# Row 35 - Hedera
$ws.Range('D35').Value = '0.0850'
$ws.Range('E35').Value = '  +0.87%  '

# Row 36 - Render Token
$ws.Range('D36').Value = '5.19'
$ws.Range('E36').Value = '  -3.13%  '

# Row 37 - First Digital USD
$ws.Range('E37').Value = '  -0.06%  '

# Row 38 - Celestia
$ws.Range('D38').Value = '18.88'
$ws.Range('E38').Value = '  +3.25%  '

# Row 39 - Lido DAO Token
$ws.Range('E39').Value = '  -2.83%  '

# Row 40 - ARBITRUM
$ws.Range('E40').Value = '  -2.53%  '

# Row 41 - Stacks
$ws.Range('E41').Value = '  +3.40%  '

# Row 42 - Stellar
$ws.Range('E42').Value = '  -1.52%  '

# Row 43 - WEMIX Token
$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  -2.26%  '

# Row 44 - now Monero (was EnergySwap; ranking order swapped with row 45)
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '119.81'
$ws.Range('E44').Value = '  -3.66%  '

# Row 45 - now EnergySwap (was Monero)
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '21.87'
$ws.Range('E45').Value = '  -5.60%  '

# Row 46 - Maker
$ws.Range('D46').Value = '2.097.80'
$ws.Range('E46').Value = '  +0.45%  '

# Row 47 - NEAR Protocol
$ws.Range('D47').Value = '3.30'
$ws.Range('E47').Value = '  -0.90%  '

# Row 48 - ApeX Protocol
$ws.Range('E48').Value = '  +1.13%  '

# Row 49 - THORChain
$ws.Range('D49').Value = '5.75'
$ws.Range('E49').Value = '  -3.46%  '

# Row 50 - SEI
$ws.Range('E50').Value = '  -2.59%  '

# Row 51 - Bitget Token
$ws.Range('E51').Value = '  +29.05%  '
